# Fix Excel file loading - skip first 4 rows, use columns directly
# Updated test file demo/data/测试数据_NC50.xlsx with 4 header rows to
# match expected format:
#   - Insert 3 new rows above the existing header row, so rows 1-4 become
#     header/info rows (row 1 keeps its position, rows 2-4 are new).
#   - Replace the old bilingual header text ("时间" / "压力") with generic
#     header labels, and give each of the 4 header rows its own "InfoN"
#     label in column B.
#   - Strip the bold/centered/bordered header styling that used to be on
#     row 1 so every header row uses plain, unformatted cells.
#   - All existing time/pressure data rows shift down by 3 rows
#     (old row 2 -> row 5 ... old row 42 -> row 45) and keep their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the current row 1 (and everything below it) down by 3 rows, opening
# up rows 1-3 for the new header rows. The newly inserted rows pick up no
# special formatting (plain/default style). The old row 1 - with its bold +
# bordered + centered styling (style index 1) - is now row 4.
$ws.Rows("1:3").Insert()

# Strip the bold/border/centered styling that rode along with the old
# header row (now at row 4) so every header row ends up plain/unformatted.
$ws.Range("A4:B4").ClearFormats()

# Row 1: new plain header text (replaces the old "时间" / "压力" values that
# now live, re-styled, in row 4 below).
$ws.Range("A1").Value = "Header Row 1"
$ws.Range("B1").Value = "Info1"

# New header rows 2-4.
$ws.Range("A2").Value = "Header Row 2"
$ws.Range("B2").Value = "Info2"

$ws.Range("A3").Value = "Header Row 3"
$ws.Range("B3").Value = "Info3"

$ws.Range("A4").Value = "Header Row 4"
$ws.Range("B4").Value = "Info4"
